{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Delete the \"B\u00e1o c\u00e1o tu\u1ea7n 2/3/4\" paragraphs entirely (including their\n// paragraph marks), leaving \"tu\u1ea7n 1\" and the trailing empty paragraph intact.\nconst targets = [\n  \"B\u00e1o c\u00e1o tu\u1ea7n 2: load dataset v\u00e0 hi\u1ec3n th\u1ecb tr\u00ean ng\u00f4n ng\u1eef l\u1eadp tr\u00ecnh\",\n  \"B\u00e1o c\u00e1o tu\u1ea7n 3: test l\u1ea7n 1\",\n  \"B\u00e1o c\u00e1o tu\u1ea7n 4: test l\u1ea7n 2 \u0111\u1ec3 hi\u1ec3u s\u00e2u\",\n];\n\nfor (const p of paragraphs.items) {\n  if (targets.includes(p.text)) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Delete the \"B\u00e1o c\u00e1o tu\u1ea7n 2/3/4\" paragraphs entirely (including their\n# paragraph marks), leaving \"tu\u1ea7n 1\" and the trailing empty paragraph intact.\n$targets = @(\n    \"B\u00e1o c\u00e1o tu\u1ea7n 2: load dataset v\u00e0 hi\u1ec3n th\u1ecb tr\u00ean ng\u00f4n ng\u1eef l\u1eadp tr\u00ecnh\",\n    \"B\u00e1o c\u00e1o tu\u1ea7n 3: test l\u1ea7n 1\",\n    \"B\u00e1o c\u00e1o tu\u1ea7n 4: test l\u1ea7n 2 \u0111\u1ec3 hi\u1ec3u s\u00e2u\"\n)\n\n$i = $d.Paragraphs.Count\nwhile ($i -ge 1) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($targets -contains $text) {\n        $p.Range.Delete()\n    }\n    $i = $i - 1\n}\n"}
